$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# List of (row, value) pairs for the new "2021" column O, mirroring the
# formatting of the existing 2020 column (N) in each row.
$rows = @(
    @{ Row = 3;  Value = 2021 },
    @{ Row = 5;  Value = 2148.1999999999998 },
    @{ Row = 6;  Value = 109.5 },
    @{ Row = 7;  Value = 210.1 },
    @{ Row = 8;  Value = 196 },
    @{ Row = 9;  Value = 209 },
    @{ Row = 10; Value = 300.2 },
    @{ Row = 11; Value = 302.89999999999998 },
    @{ Row = 12; Value = 786 },
    @{ Row = 13; Value = 27.7 },
    @{ Row = 14; Value = 6.8 },
    @{ Row = 16; Value = 26.9 },
    @{ Row = 17; Value = 15.9 },
    @{ Row = 18; Value = 21.7 },
    @{ Row = 19; Value = 29.9 },
    @{ Row = 20; Value = 30.2 },
    @{ Row = 21; Value = 24 },
    @{ Row = 22; Value = 31.6 },
    @{ Row = 23; Value = 30.3 },
    @{ Row = 24; Value = 20.7 },
    @{ Row = 25; Value = 12 }
)

foreach ($item in $rows) {
    $r = $item.Row
    $srcCell = $ws.Cells.Item($r, 14)   # column N
    $dstCell = $ws.Cells.Item($r, 15)   # column O
    $srcCell.Copy()
    $dstCell.PasteSpecial($xlPasteFormats)
    $dstCell.Value = $item.Value
}

# Row 15 has an empty O cell, but still formatted like N15.
$ws.Cells.Item(15, 14).Copy()
$ws.Cells.Item(15, 15).PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false

# Update selection to match the recorded cursor position.
$ws.Range("Q20").Select()
